# PBL.xlsx — [12] Berechnung der akkumulierten Aufwandschätzung bis zur
# Fertigstellung jedes PBIs.
#
# The PBI in row 8 (Id 12) is reworded from "alle PBIs" to "jedes PBI".
# Editing the cell's text in place causes Excel to re-home that shared
# string at the end of the shared-strings table, which is why the diff
# also shows index shifts on the untouched neighbouring rows (6-8) -
# setting the .Value below reproduces that naturally.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = "Als PO möchte ich für jedes PBI einen Forecast bis zu welchem Sprint dieses fertiggestellt wird."

# Reflect the author's final scroll/selection position (row 9 in view).
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B9").Select()
